$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 27 with the PRIME crypto asset entry (Gaming category)
$ws.Range("A27").Value = "PRIME"
$ws.Range("B27").Value = 68.66
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = "Gaming"

# Update the active selection to match the target (single cell D27)
$ws.Range("D27").Select()
